$wb = $excel.ActiveWorkbook

# The edited sheet is the 3rd sheet in the workbook ("2010-18"), which is
# also the active/tab-selected sheet.
$ws = $wb.Worksheets.Item(3)

# --- New row of data (row 32) -------------------------------------------
# Column A: model name, Column B: simulation run name, Column C: weather years
$ws.Range("A32").Value = "CW3M C787+"
$ws.Range("B32").Value = "Baseline 2010-18 "
$ws.Range("C32").Value = "2010-18"

$ws.Range("D32").Value = 493.80679655555559
$ws.Range("E32").Value = 2094.2995878888887
$ws.Range("F32").Value = 6.3996251111111109
$ws.Range("G32").Value = 332.04267011111119
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0.24795422222222227
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 596.39469400000007
$ws.Range("L32").Value = 92.091200777777772
$ws.Range("M32").Value = 1741.6711831111111
$ws.Range("N32").Value = 495.68789333333342
$ws.Range("O32").Value = 16151.769531111109
$ws.Range("P32").Value = 2215.2681748888886
$ws.Range("Q32").Value = -0.95166244444444414
$ws.Range("R32").Value = -0.0003254444444444445

# Match number formats / highlight fill used for row 32, following the
# same per-column style pattern already used in the sheet (style 3/5/8/10).
$ws.Range("E32:L32").NumberFormat = "0.00"
$ws.Range("D32").NumberFormat = "0.00"
$ws.Range("I32").NumberFormat = "0.00"
$ws.Range("M32").NumberFormat = "0.00"
$ws.Range("N32").NumberFormat = "0.00"
$ws.Range("Q32").NumberFormat = "0.00"
$ws.Range("O32:P32").NumberFormat = "0"
$ws.Range("R32").NumberFormat = "0.000000"

$yellow = 65535
$ws.Range("D32").Interior.Color = $yellow
$ws.Range("I32").Interior.Color = $yellow
$ws.Range("M32").Interior.Color = $yellow
$ws.Range("N32").Interior.Color = $yellow
$ws.Range("Q32").Interior.Color = $yellow
$ws.Range("R32").Interior.Color = $yellow

# --- Column width for column A -------------------------------------------
# Target OOXML width is 12.6640625; the closest value this runtime can
# produce via ColumnWidth (quantized in 1/6 character-width steps) is
# 12.666666666666666, reached with a ColumnWidth of 11.8.
$ws.Range("A1").EntireColumn.ColumnWidth = 11.8

# --- Selection / active cell ---------------------------------------------
# Reflects the user having extended the selection from R32 to include Q32.
$ws.Range("Q32:R32").Select() | Out-Null
